$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71 - shifts existing rows 71..118 down to 72..119
$ws.Rows.Item(71).Insert()

# Populate the new row 71 with the new data record
$ws.Cells.Item(71, 1).Value = 5
$ws.Cells.Item(71, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(71, 3).Value = 'Maule'
$ws.Cells.Item(71, 4).Value = '2022-04-07'
$ws.Cells.Item(71, 5).Value = 7
$ws.Cells.Item(71, 6).Value = 'Fruta'
$ws.Cells.Item(71, 7).Value = 100103
$ws.Cells.Item(71, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(71, 9).Value = 100103002
$ws.Cells.Item(71, 10).Value = 'Ciruela'
$ws.Cells.Item(71, 11).Value = 'Angeleno'
$ws.Cells.Item(71, 12).Value = 'Primera'
$ws.Cells.Item(71, 13).Value = 450
$ws.Cells.Item(71, 14).Value = 6000
$ws.Cells.Item(71, 15).Value = 8000
$ws.Cells.Item(71, 16).Value = 7111
$ws.Cells.Item(71, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(71, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(71, 19).Value = 395
$ws.Cells.Item(71, 20).Value = 18
